$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Match the width of the column immediately to the left (M) so the newly
# inserted column N keeps a sensible width (mirrors Excel's own behaviour
# when a column is inserted next to an existing, explicitly-sized column).
$mWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column before column N ("Late"), pushing the existing
# N/O/P ("Late", heading, Outstanding) columns one to the right.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# The repayment schedule tab becomes the active/selected sheet, with the
# cursor left on J16.
$ws.Activate()
$ws.Range("J16").Select() | Out-Null
